$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.535.60'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '1.875.68'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4870'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.77%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2894'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06671'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.05%  '
$ws.Range("D10").Value = '1.874.12'
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.62'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07224'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '89.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.002'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6547'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.79%  '
$ws.Range("D16").Value = '30.466.86'
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007823'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("D20").Value = '2.114.52'
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '211.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +18.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.734'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.143'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.379'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.830'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.415'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.262'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09051'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.928'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.68%  '
$ws.Range("E33").Value = '  -1.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7265'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.078'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.689'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("E37").Value = '  -2.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.663'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9174'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.038'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4422'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.739'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9943'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1327'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.339'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.98%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4019'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.31%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05838'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.595'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.412'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.44%  '
